$wb = $excel.ActiveWorkbook

# Row 8 on ALC  (@@ -1024,22 +1024,22 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 18
$ws.Range("I8").Value = 18
$ws.Range("K8").Value = 54
$ws.Range("M8").Value = 85

# Row 21 on ALC  (@@ -1655,22 +1655,25 @@)
$ws.Range("H21").Value = 4253.75
$ws.Range("I21").Value = 2338.6667
$ws.Range("J21").Value = 9999
$ws.Range("K21").Value = 2338.6667
$ws.Range("L21").Value = 9999
$ws.Range("M21").Value = -1870.6667
$ws.Range("N21").Value = -10935

# Row 23 on ALC  (@@ -1750,22 +1753,25 @@)
$ws.Range("H23").Value = 4253.75
$ws.Range("I23").Value = 2338.6667
$ws.Range("J23").Value = 9999
$ws.Range("K23").Value = 2338.6667
$ws.Range("L23").Value = 9999
$ws.Range("M23").Value = -2104.6667
$ws.Range("N23").Value = -10467

# Row 31 on ALC  (@@ -2139,22 +2145,22 @@)
$ws.Range("H31").Value = 932
$ws.Range("I31").Value = 932
$ws.Range("K31").Value = 2796
$ws.Range("M31").Value = -2566

# Row 51 on ALC  (@@ -3146,25 +3152,25 @@)
$ws.Range("H51").Value = 17999.666
$ws.Range("J51").Value = 24999.5
$ws.Range("L51").Value = 24999.5
$ws.Range("N51").Value = -25967.5

# Row 86 on ALC  (@@ -4891,25 +4897,25 @@)
$ws.Range("H86").Value = 1600
$ws.Range("J86").Value = 1700
$ws.Range("L86").Value = 1700
$ws.Range("N86").Value = -3946

# Row 89 on ALC  (@@ -5044,25 +5050,25 @@)
$ws.Range("H89").Value = 1600
$ws.Range("J89").Value = 1700
$ws.Range("L89").Value = 8500
$ws.Range("N89").Value = -19732

# Row 98 on ALC  (@@ -5494,25 +5500,19 @@)
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("M98").ClearContents()
$ws.Range("N98").ClearContents()

# Row 107 on ALC  (@@ -5950,25 +5950,25 @@)
$ws.Range("H107").Value = 1390.8
$ws.Range("I107").Value = 1198
$ws.Range("J107").Value = 1969.2
$ws.Range("K107").Value = 1198
$ws.Range("L107").Value = 1969.2
$ws.Range("M107").Value = 722
$ws.Range("N107").Value = -5809.2

# Row 122 on ALC  (@@ -6706,25 +6706,19 @@)
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()

# Row 2 on ARM  (@@ -7804,25 +7798,25 @@)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3437.375
$ws.Range("I2").Value = 2249.6667
$ws.Range("J2").Value = 7000.5
$ws.Range("K2").Value = 2249.6667
$ws.Range("L2").Value = 7000.5
$ws.Range("M2").Value = -2136.6667
$ws.Range("N2").Value = -7226.5

# Row 45 on ARM  (@@ -9923,25 +9917,25 @@)
$ws.Range("H45").Value = 3811.75
$ws.Range("I45").Value = 3249.2222
$ws.Range("J45").Value = 5499.3335
$ws.Range("K45").Value = 3249.2222
$ws.Range("L45").Value = 5499.3335
$ws.Range("M45").Value = -2872.2222
$ws.Range("N45").Value = -6253.3335

# Row 61 on ARM  (@@ -10704,25 +10698,25 @@)
$ws.Range("H61").Value = 2093.0977
$ws.Range("I61").Value = 2041.7059
$ws.Range("J61").Value = 2342.7144
$ws.Range("K61").Value = 2041.7059
$ws.Range("L61").Value = 2342.7144
$ws.Range("M61").Value = -1829.7059
$ws.Range("N61").Value = -2766.7144

# Row 76 on ARM  (@@ -11427,22 +11421,22 @@)
$ws.Range("H76").Value = 71500
$ws.Range("J76").Value = 71500
$ws.Range("L76").Value = 71500
$ws.Range("N76").Value = -72176

# Row 79 on ARM  (@@ -11577,22 +11571,22 @@)
$ws.Range("H79").Value = 71500
$ws.Range("J79").Value = 71500
$ws.Range("L79").Value = 71500
$ws.Range("N79").Value = -73840

# Row 116 on ARM  (@@ -13366,25 +13360,25 @@)
$ws.Range("H116").Value = 3437.375
$ws.Range("I116").Value = 2249.6667
$ws.Range("J116").Value = 7000.5
$ws.Range("K116").Value = 2249.6667
$ws.Range("L116").Value = 7000.5
$ws.Range("M116").Value = 44.33329999999978
$ws.Range("N116").Value = -11588.5

# Row 122 on ARM  (@@ -13651,22 +13645,22 @@)
$ws.Range("H122").Value = 2805.1177
$ws.Range("I122").Value = 2650.1428
$ws.Range("K122").Value = 7950.428400000001
$ws.Range("M122").Value = -5500.428400000001

# Row 136 on ARM  (@@ -14334,25 +14328,25 @@)
$ws.Range("H136").Value = 2093.0977
$ws.Range("I136").Value = 2041.7059
$ws.Range("J136").Value = 2342.7144
$ws.Range("K136").Value = 6125.1177
$ws.Range("L136").Value = 7028.1432
$ws.Range("M136").Value = -3575.1177
$ws.Range("N136").Value = -12128.1432

# Row 3 on BSM  (@@ -14768,25 +14762,25 @@)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3437.375
$ws.Range("I3").Value = 2249.6667
$ws.Range("J3").Value = 7000.5
$ws.Range("K3").Value = 2249.6667
$ws.Range("L3").Value = 7000.5
$ws.Range("M3").Value = -2135.6667
$ws.Range("N3").Value = -7228.5

# Row 75 on BSM  (@@ -18281,25 +18275,25 @@)
$ws.Range("H75").Value = 56610.875
$ws.Range("I75").Value = 7662.6665
$ws.Range("J75").Value = 85979.8
$ws.Range("K75").Value = 7662.6665
$ws.Range("L75").Value = 85979.8
$ws.Range("M75").Value = -6726.6665
$ws.Range("N75").Value = -87851.8

# Row 78 on BSM  (@@ -18431,25 +18425,25 @@)
$ws.Range("H78").Value = 56610.875
$ws.Range("I78").Value = 7662.6665
$ws.Range("J78").Value = 85979.8
$ws.Range("K78").Value = 22987.9995
$ws.Range("L78").Value = 257939.4
$ws.Range("M78").Value = -18307.9995
$ws.Range("N78").Value = -267299.4

# Row 80 on BSM  (@@ -18532,25 +18526,25 @@)
$ws.Range("H80").Value = 598
$ws.Range("I80").Value = 997.3333
$ws.Range("J80").Value = 478.2
$ws.Range("K80").Value = 997.3333
$ws.Range("L80").Value = 478.2
$ws.Range("M80").Value = 0.6666999999999916
$ws.Range("N80").Value = -2474.2

# Row 83 on BSM  (@@ -18685,25 +18679,25 @@)
$ws.Range("H83").Value = 598
$ws.Range("I83").Value = 997.3333
$ws.Range("J83").Value = 478.2
$ws.Range("K83").Value = 4986.6665
$ws.Range("L83").Value = 2391
$ws.Range("M83").Value = 5.333499999999731
$ws.Range("N83").Value = -12375

# Row 107 on BSM  (@@ -19885,25 +19879,25 @@)
$ws.Range("H107").Value = 1673.6666
$ws.Range("I107").Value = 1484.2
$ws.Range("J107").Value = 2265.75
$ws.Range("K107").Value = 1484.2
$ws.Range("L107").Value = 2265.75
$ws.Range("M107").Value = 435.8
$ws.Range("N107").Value = -6105.75

# Row 134 on BSM  (@@ -21175,25 +21169,25 @@)
$ws.Range("H134").Value = 4369
$ws.Range("I134").Value = 2112.8333
$ws.Range("J134").Value = 7753.25
$ws.Range("K134").Value = 6338.499899999999
$ws.Range("L134").Value = 23259.75
$ws.Range("M134").Value = -3803.499899999999
$ws.Range("N134").Value = -28329.75

# Row 58 on CRP  (@@ -24408,25 +24402,25 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 11959185
$ws.Range("I58").Value = 2507.2942
$ws.Range("J58").Value = 45836440
$ws.Range("K58").Value = 2507.2942
$ws.Range("L58").Value = 45836440
$ws.Range("M58").Value = -2304.2942
$ws.Range("N58").Value = -45836846

# Row 62 on CRP  (@@ -24607,22 +24601,22 @@)
$ws.Range("H62").Value = 3813.7144
$ws.Range("I62").Value = 3678
$ws.Range("K62").Value = 3678
$ws.Range("M62").Value = -3054

# Row 65 on CRP  (@@ -24751,22 +24745,22 @@)
$ws.Range("H65").Value = 3813.7144
$ws.Range("I65").Value = 3678
$ws.Range("K65").Value = 18390
$ws.Range("M65").Value = -15270

# Row 94 on CRP  (@@ -26151,25 +26145,25 @@)
$ws.Range("H94").Value = 2560.6667
$ws.Range("J94").Value = 2577.7693
$ws.Range("L94").Value = 2577.7693
$ws.Range("N94").Value = -3479.7693

# Row 134 on CRP  (@@ -28129,25 +28123,25 @@)
$ws.Range("H134").Value = 2833.3333
$ws.Range("I134").Value = 1000
$ws.Range("J134").Value = 3750
$ws.Range("K134").Value = 3000
$ws.Range("L134").Value = 11250
$ws.Range("M134").Value = -465
$ws.Range("N134").Value = -16320

# Row 136 on CRP  (@@ -28227,25 +28221,25 @@)
$ws.Range("H136").Value = 11959185
$ws.Range("I136").Value = 2507.2942
$ws.Range("J136").Value = 45836440
$ws.Range("K136").Value = 7521.882599999999
$ws.Range("L136").Value = 137509320
$ws.Range("M136").Value = -4971.882599999999
$ws.Range("N136").Value = -137514420

# Row 37 on CUL  (@@ -30393,22 +30387,22 @@)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 119511.43
$ws.Range("J37").Value = 119511.43
$ws.Range("L37").Value = 358534.29
$ws.Range("N37").Value = -358758.29

# Row 39 on CUL  (@@ -30494,25 +30488,25 @@)
$ws.Range("H39").Value = 2777.7778
$ws.Range("I39").Value = 2000
$ws.Range("J39").Value = 3000
$ws.Range("K39").Value = 6000
$ws.Range("L39").Value = 9000
$ws.Range("M39").Value = -5706
$ws.Range("N39").Value = -9588

# Row 40 on CUL  (@@ -30546,25 +30540,25 @@)
$ws.Range("H40").Value = 133.41176
$ws.Range("I40").Value = 50.666668
$ws.Range("J40").Value = 226.5
$ws.Range("K40").Value = 202.666672
$ws.Range("L40").Value = 906
$ws.Range("M40").Value = -133.666672
$ws.Range("N40").Value = -1044

# Row 128 on CUL  (@@ -34957,22 +34951,22 @@)
$ws.Range("H128").Value = 159890
$ws.Range("I128").Value = 159890
$ws.Range("K128").Value = 479670
$ws.Range("M128").Value = -474690

# Row 59 on GSM  (@@ -38626,25 +38620,25 @@)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H59").Value = 17550
$ws.Range("J59").Value = 25000
$ws.Range("L59").Value = 25000
$ws.Range("N59").Value = -26166

# Row 70 on GSM  (@@ -39165,25 +39159,25 @@)
$ws.Range("H70").Value = 7441.2354
$ws.Range("I70").Value = 4864.364
$ws.Range("J70").Value = 12165.5
$ws.Range("K70").Value = 4864.364
$ws.Range("L70").Value = 12165.5
$ws.Range("M70").Value = -4594.364
$ws.Range("N70").Value = -12705.5

# Row 73 on GSM  (@@ -39312,25 +39306,25 @@)
$ws.Range("H73").Value = 7441.2354
$ws.Range("I73").Value = 4864.364
$ws.Range("J73").Value = 12165.5
$ws.Range("K73").Value = 4864.364
$ws.Range("L73").Value = 12165.5
$ws.Range("M73").Value = -3928.364
$ws.Range("N73").Value = -14037.5

# Row 80 on GSM  (@@ -39646,19 +39640,22 @@)
$ws.Range("H80").Value = 7999
$ws.Range("I80").Value = 7999
$ws.Range("K80").Value = 7999
$ws.Range("M80").Value = -7001

# Row 83 on GSM  (@@ -39787,19 +39784,22 @@)
$ws.Range("H83").Value = 7999
$ws.Range("I83").Value = 7999
$ws.Range("K83").Value = 39995
$ws.Range("M83").Value = -35003

# Row 132 on GSM  (@@ -42179,25 +42179,25 @@)
$ws.Range("H132").Value = 2240.4827
$ws.Range("I132").Value = 1652.3043
$ws.Range("J132").Value = 4495.1665
$ws.Range("K132").Value = 4956.9129
$ws.Range("L132").Value = 13485.4995
$ws.Range("M132").Value = -2426.9129
$ws.Range("N132").Value = -18545.4995

# Row 40 on LTW  (@@ -44610,22 +44610,22 @@)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6262.6665
$ws.Range("I40").Value = 4398
$ws.Range("K40").Value = 4398
$ws.Range("M40").Value = -4262

# Row 55 on LTW  (@@ -45360,25 +45360,25 @@)
$ws.Range("H55").Value = 819.7778
$ws.Range("I55").Value = 718.4
$ws.Range("J55").Value = 946.5
$ws.Range("K55").Value = 718.4
$ws.Range("L55").Value = 946.5
$ws.Range("M55").Value = -545.4
$ws.Range("N55").Value = -1292.5

# Row 122 on LTW  (@@ -48661,22 +48661,22 @@)
$ws.Range("H122").Value = 4115.359
$ws.Range("I122").Value = 3757.6428
$ws.Range("K122").Value = 11272.9284
$ws.Range("M122").Value = -8822.928400000001

# Row 136 on LTW  (@@ -49347,22 +49347,22 @@)
$ws.Range("H136").Value = 5171
$ws.Range("I136").Value = 2776.9473
$ws.Range("K136").Value = 8330.841899999999
$ws.Range("M136").Value = -5780.841899999999

# Row 126 on WVR  (@@ -55859,26 +55859,23 @@)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 55576056
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()

# Row 132 on WVR  (@@ -56150,22 +56147,22 @@)
$ws.Range("H132").Value = 5863.273
$ws.Range("I132").Value = 5849.6
$ws.Range("K132").Value = 17548.8
$ws.Range("M132").Value = -15018.8

# Row 133 on WVR  (@@ -56202,22 +56199,22 @@)
$ws.Range("H133").Value = 79650.75
$ws.Range("J133").Value = 79650.75
$ws.Range("L133").Value = 79650.75
$ws.Range("N133").Value = -89770.75
